# Apply the cryptos-list refresh described in the commit:
# "Updated cryptos list on Fri Mar  3 08:49:42 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Forces the cell to keep its original plain-text representation
    # (e.g. "291.16", "0.00001139") instead of being auto-coerced into a
    # floating point number, then restores the default (unstyled) cell
    # formatting so no stray style index is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '22.404.48'
$ws.Range("E2").Value = '  -4.52%  '
$ws.Range("D3").Value = '1.571.24'
$ws.Range("E3").Value = '  -4.61%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  -0.05%  '
Set-TextValue $ws.Range("D6") '291.16'
$ws.Range("E6").Value = '  -2.52%  '
Set-TextValue $ws.Range("D7") '0.3653'
$ws.Range("E7").Value = '  -3.46%  '
Set-TextValue $ws.Range("D8") '49.38'
$ws.Range("E8").Value = '  -0.82%  '
Set-TextValue $ws.Range("D9") '0.3378'
$ws.Range("E9").Value = '  -5.01%  '
Set-TextValue $ws.Range("D10") '1.172'
$ws.Range("E10").Value = '  -3.90%  '
Set-TextValue $ws.Range("D11") '0.07588'
$ws.Range("E11").Value = '  -6.27%  '
$ws.Range("E12").Value = '  +0.07%  '
Set-TextValue $ws.Range("D13") '21.17'
$ws.Range("E13").Value = '  -4.07%  '
Set-TextValue $ws.Range("D14") '6.060'
$ws.Range("E14").Value = '  -5.26%  '
Set-TextValue $ws.Range("D15") '6.875'
$ws.Range("E15").Value = '  -6.35%  '
Set-TextValue $ws.Range("D16") '0.00001139'
$ws.Range("E16").Value = '  -4.80%  '
$ws.Range("D17").Value = '1.569.58'
$ws.Range("E17").Value = '  -5.12%  '
Set-TextValue $ws.Range("D18") '89.04'
$ws.Range("E18").Value = '  -8.52%  '
Set-TextValue $ws.Range("D19") '0.06739'
$ws.Range("E19").Value = '  -3.03%  '
$ws.Range("E20").Value = '  -0.04%  '
Set-TextValue $ws.Range("D21") '6.270'
$ws.Range("E21").Value = '  -7.25%  '
Set-TextValue $ws.Range("D22") '16.46'
$ws.Range("E22").Value = '  -4.96%  '
Set-TextValue $ws.Range("D23") '0.5242'
$ws.Range("E23").Value = '  -8.79%  '
Set-TextValue $ws.Range("D24") '12.01'
$ws.Range("E24").Value = '  -3.43%  '
$ws.Range("D25").Value = '22.415.88'
$ws.Range("E25").Value = '  -4.54%  '
$ws.Range("E26").Value = '  -4.18%  '
Set-TextValue $ws.Range("D27") '3.002'
$ws.Range("E27").Value = '  +3.38%  '
Set-TextValue $ws.Range("D28") '19.89'
$ws.Range("E28").Value = '  -4.89%  '
Set-TextValue $ws.Range("D29") '144.30'
$ws.Range("E29").Value = '  -5.70%  '
Set-TextValue $ws.Range("D30") '4.993'
$ws.Range("E30").Value = '  -4.22%  '
Set-TextValue $ws.Range("D31") '125.20'
$ws.Range("E31").Value = '  -5.80%  '
$ws.Range("D32").Value = '1.745.10'
$ws.Range("E32").Value = '  -4.75%  '
Set-TextValue $ws.Range("D33") '1.048'
$ws.Range("E33").Value = '  +4.90%  '
Set-TextValue $ws.Range("D34") '6.304'
$ws.Range("E34").Value = '  -8.78%  '
Set-TextValue $ws.Range("D35") '1.972'
$ws.Range("E35").Value = '  -7.05%  '
Set-TextValue $ws.Range("D36") '10.36'
$ws.Range("E36").Value = '  -9.71%  '
Set-TextValue $ws.Range("D37") '0.02559'
$ws.Range("E37").Value = '  -5.80%  '
Set-TextValue $ws.Range("D38") '0.08437'
$ws.Range("E38").Value = '  -3.44%  '
Set-TextValue $ws.Range("D39") '0.2309'
$ws.Range("E39").Value = '  -4.78%  '
Set-TextValue $ws.Range("D40") '0.06533'
$ws.Range("E40").Value = '  -3.55%  '
Set-TextValue $ws.Range("D41") '5.520'
$ws.Range("E41").Value = '  -6.92%  '
Set-TextValue $ws.Range("D42") '11.86'
$ws.Range("E42").Value = '  -9.29%  '
Set-TextValue $ws.Range("D43") '1.252'
$ws.Range("E43").Value = '  -3.99%  '
Set-TextValue $ws.Range("D44") '0.6397'
$ws.Range("E44").Value = '  -7.04%  '
$ws.Range("E45").Value = '  -6.48%  '
$ws.Range("E46").Value = '  -0.08%  '
Set-TextValue $ws.Range("D47") '0.6030'
$ws.Range("E47").Value = '  -5.20%  '
Set-TextValue $ws.Range("D48") '3.774'
$ws.Range("E48").Value = '  -3.50%  '
Set-TextValue $ws.Range("D49") '2.139'
$ws.Range("E49").Value = '  -5.10%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D50") '122.45'
$ws.Range("E50").Value = '  -3.74%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue $ws.Range("D51") '1.209'
$ws.Range("E51").Value = '  +2.55%  '
